$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 116, shifting rows 116-132 down to 117-133
$ws.Rows.Item(116).Insert()

# Copy number format from the row below (row 117, which is the old row 116) for the date cell
$ws.Cells.Item(116, 4).NumberFormat = $ws.Cells.Item(117, 4).NumberFormat

# Populate new row 116 with data
$ws.Cells.Item(116, 1).Value = 10
$ws.Cells.Item(116, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(116, 3).Value = "La Araucanía"
$ws.Cells.Item(116, 4).Value = 44474
$ws.Cells.Item(116, 5).Value = 9
$ws.Cells.Item(116, 6).Value = "Fruta"
$ws.Cells.Item(116, 7).Value = 100101
$ws.Cells.Item(116, 8).Value = "Berries"
$ws.Cells.Item(116, 9).Value = 100112025
$ws.Cells.Item(116, 10).Value = "Frutilla"
$ws.Cells.Item(116, 11).Value = "Sin especificar"
$ws.Cells.Item(116, 12).Value = "Primera"
$ws.Cells.Item(116, 13).Value = 300
$ws.Cells.Item(116, 14).Value = 14000
$ws.Cells.Item(116, 15).Value = 14000
$ws.Cells.Item(116, 16).Value = 14000
$ws.Cells.Item(116, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(116, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(116, 19).Value = 2000
$ws.Cells.Item(116, 20).Value = 7
